# "3- Regras de Comunicação.docx" edit
#
# Content changes applied:
#   1. The first bullet item gets an extra clause inserted just before its
#      closing period:
#        "... por telefone ou pessoalmente."
#        -> "... por telefone ou pessoalmente em horários de folga."
#   2. The second bullet item ("Reuniões presenciais entre os integrantes
#      devem ocorrer somente aos finais de semana.") is removed: its text is
#      deleted and the paragraph is demoted from the bulleted list back to a
#      plain "Normal" paragraph (so it becomes one of the blank lines at the
#      end of the document, like its neighbours).
#   3. Two of the extra blank paragraphs that used to pad the end of the
#      document are removed as a result, so the document ends with a tidy
#      set of blank paragraphs instead of a long trailing gap.

$d = $word.ActiveDocument

# 1. Extend the first bullet's sentence with the extra clause before the
#    final period.
[void]$d.Content.Find.Execute(
    "ou pessoalmente.", $true, $false, $false, $false, $false, $true, 1,
    $false, "ou pessoalmente em horários de folga.", 2)

# 2. Locate the "Reuniões presenciais..." bullet paragraph, clear its text
#    and remove its list numbering, turning it into a blank Normal
#    paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Reuniões presenciais*") {
        $paraRange = $para.Range
        $textOnly = $d.Range($paraRange.Start, $paraRange.End - 1)
        $textOnly.Text = ""
        $para.Range.ListFormat.RemoveNumbers(0)
        $para.Range.Style = "Normal"
        break
    }
}

# 3. Drop two of the now-redundant trailing blank paragraphs so the document
#    ends with three blank paragraphs after the (now unbulleted) paragraph
#    instead of four.
$lastBlank = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$lastBlank.Range.Delete()
$lastBlank2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$lastBlank2.Range.Delete()
